# Apply the author's proof-reading pass over the "Exigence des concepts"
# and "Bonus" bullet lists: a handful of singular->plural agreement fixes,
# a missing comma, and a couple of word corrections.

$d = $word.ActiveDocument

# 1. "Nous avons utilisé l'interface" -> "...utilisés l'interface"
$d.Content.Find.Execute(
    "Nous avons utilisé l’interface Incomparable",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nous avons utilisés l’interface Incomparable",
    2)

# 2. "2 Type de collection générique ont été utiliser :" -> "2 Types ... utilisés :"
$d.Content.Find.Execute(
    "2 Type de collection générique ont été utiliser :",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2 Types de collection générique ont été utilisés :",
    2)

# 2b. "Il y a plusieurs Abstract class (Evènement & Identité)" was split across
#     three runs in the source; re-assert it as a no-op replace so the engine
#     coalesces it back into a single run (matching the cleaned-up markup).
$d.Content.Find.Execute(
    "Il y a plusieurs Abstract class (Evènement & Identité)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il y a plusieurs Abstract class (Evènement & Identité)",
    2)

# 3. "héritages tous le long" -> "héritages tout au long"
$d.Content.Find.Execute(
    "héritages tous le long",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "héritages tout au long",
    2)

# 4. "Problème rencontré :" -> "Problèmes rencontrés :"
$d.Content.Find.Execute(
    "Problème rencontré :",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Problèmes rencontrés :",
    2)

# 5. "voir explication dans le diagramme" -> "voir explications dans le diagramme"
$d.Content.Find.Execute(
    "voir explication dans le diagramme",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "voir explications dans le diagramme",
    2)

# 6. "dégrader ... image afin de rendre l'interface plus sympa." ->
#    "dégradé ... images afin de rendre l'interface plus sympathique."
$d.Content.Find.Execute(
    "est un dégrader, il y a aussi un certain nombre d’image afin de rendre l’interface plus sympa.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "est un dégradé, il y a aussi un certain nombre d’images afin de rendre l’interface plus sympathique.",
    2)

# 7. "Lors de ce projet nous avons" -> "Lors de ce projet, nous avons"
$d.Content.Find.Execute(
    "Lors de ce projet nous avons",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lors de ce projet, nous avons",
    2)

# 8. "L'utilisation est fluide ... boutons dont" was split across three runs
#    in the source; re-assert it as a no-op replace to coalesce it back into
#    a single run.
$d.Content.Find.Execute(
    "L’utilisation est fluide et intuitive grâce aux nombreux boutons dont",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "L’utilisation est fluide et intuitive grâce aux nombreux boutons dont",
    2)
